$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1132309489190599
$ws.Range("H2").Value = 19.54169361336406
$ws.Range("I2").Value = 4.485680346449999
$ws.Range("G3").Value = 0.07999476026304553
$ws.Range("H3").Value = 21.03197512195841
$ws.Range("G4").Value = 0.001884135350643496
$ws.Range("H4").Value = -79.90312981825177
$ws.Range("G5").Value = 0.02139826697649108
$ws.Range("H5").Value = 259.810058802088
$ws.Range("G6").Value = -0.2184643526377449
$ws.Range("H6").Value = 1.234804979412254
$ws.Range("G7").Value = -0.2549660960028569
$ws.Range("H7").Value = -2.037173603500511
$ws.Range("G8").Value = -0.3874074037318599
$ws.Range("H8").Value = -4.632346272796875
$ws.Range("G9").Value = -0.4355255419159755
$ws.Range("H9").Value = -9.254388244326066
$ws.Range("G10").Value = -0.03807182547976942
$ws.Range("H10").Value = -334.9403763348403
$ws.Range("G11").Value = 0.0753273591192322
$ws.Range("H11").Value = 568.8658108017067
$ws.Range("G12").Value = 0.2179881036386403
$ws.Range("H12").Value = -4.049802361289668
$ws.Range("G13").Value = 0.2455546354235534
$ws.Range("H13").Value = -6.753972068375597
$ws.Range("G14").Value = 0.01218669362007548
$ws.Range("H14").Value = 227.3789761821235
$ws.Range("G15").Value = 0.04653306550818495
$ws.Range("H15").Value = 130.513290712286
$ws.Range("G16").Value = 0.1372015492691846
$ws.Range("H16").Value = 16.27338092267776
$ws.Range("G17").Value = 0.1769554511433631
$ws.Range("H17").Value = -19.13673263521453
$ws.Range("G18").Value = 0.06228022902889507
$ws.Range("H18").Value = 2.999383866514159
$ws.Range("G19").Value = 0.07028610898038067
$ws.Range("H19").Value = -21.98173836906189
$ws.Range("G20").Value = -0.1452702048528579
$ws.Range("H20").Value = 0.1791852800586353
$ws.Range("G21").Value = -0.1776995314925533
$ws.Range("H21").Value = 11.07460043841821
$ws.Range("G22").Value = 0.04970126971753371
$ws.Range("H22").Value = -8.61762145539535
$ws.Range("G23").Value = 0.04565979463589098
$ws.Range("H23").Value = 11.8012544308014
$ws.Range("G24").Value = 0.118877049589962
$ws.Range("H24").Value = 2.715466300940102
$ws.Range("G25").Value = 0.1417702597389913
$ws.Range("H25").Value = -6.776721440837006
$ws.Range("G26").Value = 0.05123459954186331
$ws.Range("H26").Value = -3.105686525592766
$ws.Range("G27").Value = 0.06795184545913183
$ws.Range("H27").Value = 34.64661001708708
$ws.Range("G28").Value = 0.1455671779637235
$ws.Range("H28").Value = -4.80272250838915
$ws.Range("G29").Value = 0.1914685022172364
$ws.Range("H29").Value = 12.16463892592313
$ws.Range("G30").Value = 0.003601829977364633
$ws.Range("H30").Value = -81.59146486343647
$ws.Range("G31").Value = 0.03033663655102923
$ws.Range("H31").Value = 212.5886496154944
$ws.Range("G32").Value = 0.01637555253953559
$ws.Range("H32").Value = -56.08979440407126
$ws.Range("G33").Value = 0.02496853781122375
$ws.Range("H33").Value = -4.348439946433787
$ws.Range("G34").Value = 0.1024463592681746
$ws.Range("H34").Value = -19.94377316196852
$ws.Range("G35").Value = 0.1300478565098161
$ws.Range("H35").Value = 1.078234931744494
$ws.Range("G36").Value = -0.03855189881183218
$ws.Range("H36").Value = -356.4674523054322
$ws.Range("G37").Value = 0.03075711236371647
$ws.Range("H37").Value = 100.8366949958991
$ws.Range("G38").Value = -0.0468533844501826
$ws.Range("H38").Value = -2194.020404051939
$ws.Range("G39").Value = -0.03070003599534824
$ws.Range("H39").Value = 8.109069947603134
$ws.Range("G40").Value = 0.1313343778615012
$ws.Range("H40").Value = -10.98930869121273
$ws.Range("G41").Value = 0.1486924436117502
$ws.Range("H41").Value = -7.872877941618805
$ws.Range("G42").Value = 0.04233621456857788
$ws.Range("H42").Value = -34.42844050527255
$ws.Range("G43").Value = 0.0687494686567721
$ws.Range("H43").Value = 97.78030461321106
$ws.Range("G44").Value = 0.02361918520750009
$ws.Range("H44").Value = 67.36022768116796
$ws.Range("G45").Value = 0.05176831567429117
$ws.Range("H45").Value = 26.08570108483318
$ws.Range("G46").Value = -0.03941531594621493
$ws.Range("H46").Value = 40.11638198057701
$ws.Range("G47").Value = -0.05339122331374629
$ws.Range("H47").Value = -29.24620532541438
$ws.Range("G48").Value = -0.1503394814108575
$ws.Range("H48").Value = -19.34033959453306
$ws.Range("G49").Value = -0.1372552014972804
$ws.Range("H49").Value = 30.49715217745425
$ws.Range("G50").Value = 0.1143297693833415
$ws.Range("H50").Value = 5.01142537226467
$ws.Range("G51").Value = 0.1548684740798898
$ws.Range("H51").Value = 54.45017923546123
$ws.Range("G52").Value = 0.06757346412529057
$ws.Range("H52").Value = 13.34206167379827
$ws.Range("G53").Value = 0.06220649611159036
$ws.Range("H53").Value = -7.900462785387215
$ws.Range("G54").Value = -0.1095686560659947
$ws.Range("H54").Value = -56.70444316844668
$ws.Range("G55").Value = -0.07086766981252733
$ws.Range("H55").Value = 8.245528066140922
$ws.Range("G56").Value = 0.04987099109548468
$ws.Range("H56").Value = 8.826672640264732
$ws.Range("G57").Value = 0.09004424012311593
$ws.Range("H57").Value = 1641.574503372964
